# edit.ps1
# Applies the "edits part 2 090221" change:
#   - Two occurrences of the run text "OpenMCT" (in otherwise identical
#     "AIT server and OpenMCT" textboxes) become "Open MCT".
#   - The edited run's spell-check "err" flag is cleared as a natural
#     consequence of replacing the run's text (mirrors what PowerPoint
#     itself does when a run's text content changes and is re-checked).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Update-OpenMctRun {
    param($ShapeIndex)

    $shape = $s.Shapes.Item($ShapeIndex)
    $tr = $shape.TextFrame.TextRange
    $fullText = $tr.Text
    $target = "OpenMCT"
    $startPos = $fullText.IndexOf($target)

    if ($startPos -ge 0) {
        # Characters() is 1-based.
        $oldRun = $tr.Characters($startPos + 1, $target.Length)
        $oldRun.Delete() | Out-Null

        $insertionPoint = $s.Shapes.Item($ShapeIndex).TextFrame.TextRange.Characters($startPos + 1, 0)
        $insertionPoint.InsertAfter("Open MCT") | Out-Null
    }
}

# Shape 12 ("TextBox 16" at x=2534152) and Shape 22 ("TextBox 16" at
# x=5616533) both contain the "AIT server and OpenMCT" text.
Update-OpenMctRun 12
Update-OpenMctRun 22
